$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1,1).Value = "Datos actualizados a 25 de Julio de 2020 a las 16:42"

# Row 4
$ws.Cells.Item(4,2).Value = 4264689
$ws.Cells.Item(4,3).Value = 16362
$ws.Cells.Item(4,4).Value = 2029207
$ws.Cells.Item(4,5).Value = 2086816
$ws.Cells.Item(4,7).Value = 176
$ws.Cells.Item(4,8).Value = 148666

# Row 6
$ws.Cells.Item(6,2).Value = 1364475
$ws.Cells.Item(6,3).Value = 27453
$ws.Cells.Item(6,4).Value = 870209
$ws.Cells.Item(6,5).Value = 462579
$ws.Cells.Item(6,7).Value = 281
$ws.Cells.Item(6,8).Value = 31687

# Row 16
$ws.Cells.Item(16,2).Value = 264973
$ws.Cells.Item(16,3).Value = 2201
$ws.Cells.Item(16,4).Value = 217782
$ws.Cells.Item(16,5).Value = 44488
$ws.Cells.Item(16,7).Value = 31
$ws.Cells.Item(16,8).Value = 2703

# Row 21
$ws.Cells.Item(21,2).Value = 206059
$ws.Cells.Item(21,3).Value = 99
$ws.Cells.Item(21,5).Value = 6458

# Row 23
$ws.Cells.Item(23,4).Value = 68022
$ws.Cells.Item(23,5).Value = 82651
$ws.Cells.Item(23,7).Value = 40
$ws.Cells.Item(23,8).Value = 2847

# Row 26
$ws.Cells.Item(26,2).Value = 107573
$ws.Cells.Item(26,3).Value = 2862
$ws.Cells.Item(26,4).Value = 73317
$ws.Cells.Item(26,5).Value = 29972
$ws.Cells.Item(26,7).Value = 72
$ws.Cells.Item(26,8).Value = 4284

# Row 40
$ws.Cells.Item(40,1).Value = "Republica Dominicana"
$ws.Cells.Item(40,2).Value = 60896
$ws.Cells.Item(40,3).Value = 1819
$ws.Cells.Item(40,4).Value = 27980
$ws.Cells.Item(40,5).Value = 31861
$ws.Cells.Item(40,7).Value = 19
$ws.Cells.Item(40,8).Value = 1055

# Row 41
$ws.Cells.Item(41,1).Value = "Israel"
$ws.Cells.Item(41,2).Value = 59475
$ws.Cells.Item(41,4).Value = 26797
$ws.Cells.Item(41,5).Value = 32230
$ws.Cells.Item(41,8).Value = 448

# Row 46
$ws.Cells.Item(46,4).Value = 45352
$ws.Cells.Item(46,5).Value = 4509

# Row 63
$ws.Cells.Item(63,2).Value = 22828
$ws.Cells.Item(63,3).Value = 345
$ws.Cells.Item(63,4).Value = 15608
$ws.Cells.Item(63,5).Value = 6488
$ws.Cells.Item(63,7).Value = 6
$ws.Cells.Item(63,8).Value = 732

# Row 85
$ws.Cells.Item(85,2).Value = 9097
$ws.Cells.Item(85,3).Value = 5
$ws.Cells.Item(85,5).Value = 168

# Row 92
$ws.Cells.Item(92,2).Value = 7150
$ws.Cells.Item(92,3).Value = 46
$ws.Cells.Item(92,4).Value = 5906
$ws.Cells.Item(92,5).Value = 1185
$ws.Cells.Item(92,7).Value = 1
$ws.Cells.Item(92,8).Value = 59

# Row 99
$ws.Cells.Item(99,1).Value = "Albania"
$ws.Cells.Item(99,2).Value = 4637
$ws.Cells.Item(99,3).Value = 67
$ws.Cells.Item(99,4).Value = 2637
$ws.Cells.Item(99,5).Value = 1866
$ws.Cells.Item(99,7).Value = 6
$ws.Cells.Item(99,8).Value = 134

# Row 100
$ws.Cells.Item(100,1).Value = "Republica de Africa Central"
$ws.Cells.Item(100,2).Value = 4593
$ws.Cells.Item(100,4).Value = 1483
$ws.Cells.Item(100,5).Value = 3051
$ws.Cells.Item(100,8).Value = 59

# Row 116
$ws.Cells.Item(116,1).Value = "Hong Kong"
$ws.Cells.Item(116,2).Value = 2506
$ws.Cells.Item(116,3).Value = 133
$ws.Cells.Item(116,4).Value = 1455
$ws.Cells.Item(116,5).Value = 1033
$ws.Cells.Item(116,7).Value = 2
$ws.Cells.Item(116,8).Value = 18

# Row 117
$ws.Cells.Item(117,1).Value = "Mali"
$ws.Cells.Item(117,2).Value = 2503
$ws.Cells.Item(117,4).Value = 1901
$ws.Cells.Item(117,5).Value = 479
$ws.Cells.Item(117,8).Value = 123

# Row 118
$ws.Cells.Item(118,1).Value = "Cuba"
$ws.Cells.Item(118,2).Value = 2469
$ws.Cells.Item(118,4).Value = 2341
$ws.Cells.Item(118,5).Value = 41
$ws.Cells.Item(118,8).Value = 87

# Row 119
$ws.Cells.Item(119,1).Value = "Libia"
$ws.Cells.Item(119,2).Value = 2424
$ws.Cells.Item(119,4).Value = 504
$ws.Cells.Item(119,5).Value = 1863
$ws.Cells.Item(119,7).Value = 0
$ws.Cells.Item(119,8).Value = 57

# Row 133
$ws.Cells.Item(133,1).Value = "Namibia"
$ws.Cells.Item(133,2).Value = 1687
$ws.Cells.Item(133,3).Value = 69
$ws.Cells.Item(133,4).Value = 72
$ws.Cells.Item(133,5).Value = 1608
$ws.Cells.Item(133,8).Value = 7

# Row 134
$ws.Cells.Item(134,1).Value = "Yemen"
$ws.Cells.Item(134,2).Value = 1674
$ws.Cells.Item(134,4).Value = 779
$ws.Cells.Item(134,5).Value = 426
$ws.Cells.Item(134,8).Value = 469

# Row 141
$ws.Cells.Item(141,1).Value = "Liberia"
$ws.Cells.Item(141,2).Value = 1155
$ws.Cells.Item(141,3).Value = 20
$ws.Cells.Item(141,4).Value = 631
$ws.Cells.Item(141,5).Value = 453
$ws.Cells.Item(141,8).Value = 71

# Row 142
$ws.Cells.Item(142,1).Value = "Jordania"
$ws.Cells.Item(142,2).Value = 1146
$ws.Cells.Item(142,4).Value = 1035
$ws.Cells.Item(142,5).Value = 100
$ws.Cells.Item(142,8).Value = 11

# Row 165
$ws.Cells.Item(165,2).Value = 348
$ws.Cells.Item(165,3).Value = 2
$ws.Cells.Item(165,4).Value = 288

# Row 181
$ws.Cells.Item(181,2).Value = 147
$ws.Cells.Item(181,3).Value = 5
$ws.Cells.Item(181,5).Value = 11

# Row 198
$ws.Cells.Item(198,1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(198,2).Value = 39
$ws.Cells.Item(198,3).Value = 7
$ws.Cells.Item(198,4).Value = 11
$ws.Cells.Item(198,5).Value = 28
$ws.Cells.Item(198,8).Value = 0

# Row 199
$ws.Cells.Item(199,1).Value = "Guam"
$ws.Cells.Item(199,4).Value = 0
$ws.Cells.Item(199,5).Value = 31
$ws.Cells.Item(199,8).Value = 1

# Row 210
$ws.Cells.Item(210,1).Value = "Islas Malvinas"

# Row 211
$ws.Cells.Item(211,1).Value = "Groenlandia"
